# Insert a new data row at row 381 (shifting existing rows 381-454 down to 382-455)
# and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(381).Insert()

$row = $ws.Rows.Item(381)

$row.Cells.Item(1, 1).Value = 9
$row.Cells.Item(1, 2).Value = "Vega Central Mapocho de Santiago"
$row.Cells.Item(1, 3).Value = "Metropolitana"
$row.Cells.Item(1, 4).Value = 44641
$row.Cells.Item(1, 5).Value = 13
$row.Cells.Item(1, 6).Value = 100114014
$row.Cells.Item(1, 7).Value = "Betarraga"
$row.Cells.Item(1, 8).Value = "Sin especificar"
$row.Cells.Item(1, 9).Value = "Primera"
$row.Cells.Item(1, 10).Value = 4300
$row.Cells.Item(1, 11).Value = 110
$row.Cells.Item(1, 12).Value = 120
$row.Cells.Item(1, 13).Value = 115
$row.Cells.Item(1, 14).Value = "`$/unidad"
$row.Cells.Item(1, 15).Value = "Región Metropolitana"
$row.Cells.Item(1, 16).Value = 115
$row.Cells.Item(1, 17).Value = 1
$row.Cells.Item(1, 18).Value = "Hortaliza"
